$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Total" header in column T (row 1)
$ws.Range("T1").Value = "Total"

# New row 7: "Outros" category
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 6930
$ws.Range("C7").Value = 320
$ws.Range("D7").Value = 506
$ws.Range("E7").Value = 2327
$ws.Range("F7").Value = 2993
$ws.Range("G7").Value = 3004
$ws.Range("H7").Value = 3346
$ws.Range("I7").Value = 3694
$ws.Range("J7").Value = 3973
$ws.Range("K7").Value = 4490
$ws.Range("L7").Value = 5183
$ws.Range("M7").Value = 5514
$ws.Range("N7").Value = 5490
$ws.Range("O7").Value = 5595
$ws.Range("P7").Value = 5619
$ws.Range("Q7").Value = 6618
$ws.Range("R7").Value = 20538
$ws.Range("S7").Value = 547
$ws.Range("T7").Value = 86687

# New row 8: "Total" row (column sums)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 7736
$ws.Range("C8").Value = 516
$ws.Range("D8").Value = 778
$ws.Range("E8").Value = 2955
$ws.Range("F8").Value = 3852
$ws.Range("G8").Value = 4131
$ws.Range("H8").Value = 5128
$ws.Range("I8").Value = 6454
$ws.Range("J8").Value = 8376
$ws.Range("K8").Value = 11385
$ws.Range("L8").Value = 15925
$ws.Range("M8").Value = 20670
$ws.Range("N8").Value = 24411
$ws.Range("O8").Value = 27073
$ws.Range("P8").Value = 28408
$ws.Range("Q8").Value = 32367
$ws.Range("R8").Value = 86727
$ws.Range("S8").Value = 753
$ws.Range("T8").Value = 287645

# Totals for existing rows 2-6 (column T)
$ws.Range("T2").Value = 85071
$ws.Range("T3").Value = 10039
$ws.Range("T4").Value = 39330
$ws.Range("T5").Value = 13628
$ws.Range("T6").Value = 52890
